# LOM3050.docx edit: rotate the body text of several sections into new
# homes (per the commit's unified diff) and restructure the "Critério"
# and "Bibliografia" paragraphs.
#
# Strategy: address each target paragraph by its (stable) index, find the
# byte offset of the run(s) we need to touch by locating an anchor
# substring inside the paragraph's own Range.Text, then assign .Text on a
# narrow $d.Range(start, end) so only that run's wording changes. Runs
# inside one paragraph are edited back-to-front so earlier offsets stay
# valid after a length-changing edit. (Arithmetic is always pre-computed
# into a plain variable before being handed to $d.Range(...) / .Text=.)

$d = $word.ActiveDocument
$vtab = [char]11   # how Word represents a <w:br/> inside Range.Text

if ($d.Paragraphs.Count -ne 16) {
    Write-Output "WARNING: unexpected paragraph count $($d.Paragraphs.Count)"
}

# --- "Objetivos" body (paragraph 6): single run, whole-paragraph swap ---
$pObjetivos = $d.Paragraphs.Item(6).Range
$s = $pObjetivos.Start
$e = $pObjetivos.End - 1
$r = $d.Range($s, $e)
$r.Text = "A definir, de acordo com o tópico programado."

# --- "Docente(s) Responsável(eis)" body (paragraph 8): two runs ---
# run1 + <w:br/>, then run2 (no trailing break)
$pDocentes = $d.Paragraphs.Item(8).Range
$docentesText = $pDocentes.Text
$docentesBr = $docentesText.IndexOf($vtab)
if ($docentesBr -lt 0) {
    Write-Output "WARNING: break not found in Docentes paragraph"
}

# edit run2 first (tail of the paragraph) so run1's offsets don't move
$s = $pDocentes.Start + $docentesBr + 1
$e = $pDocentes.End - 1
$r = $d.Range($s, $e)
$r.Text = "O conteúdo desta disciplina (optativa)será de acordo com o tópico a ser programado, devendo abordar assuntos complementares ao conteúdo regular do curso de graduação."

$s = $pDocentes.Start
$e = $pDocentes.Start + $docentesBr
$r = $d.Range($s, $e)
$r.Text = "Complementar a formação dos alunos em Engenharia de Materiais abordando, com maior profundidade, tópicos atuais e relevantes e atualizar com temas no estado da arte."

# --- "Programa resumido" body (paragraph 10): single run ---
$pResumido = $d.Paragraphs.Item(10).Range
$s = $pResumido.Start
$e = $pResumido.End - 1
$r = $d.Range($s, $e)
$r.Text = "Este curso deverá conter avaliações escritas e desenvolvimento de Estudo de Casos ou Projetos na área de Engenharia de Materiais. Sendo necessário aplicar pelo menos dois tipos de avaliações diferentes."

# --- "Programa" body (paragraph 12): single run ---
$pPrograma = $d.Paragraphs.Item(12).Range
$s = $pPrograma.Start
$e = $pPrograma.End - 1
$r = $d.Range($s, $e)
$r.Text = "A média do semestre será computada com base na relação: M=(A1+A2)/2"

# --- "Avaliação" body (paragraph 14): six runs --
#   [bold]"Método: "  [plain]value+br  [bold]"Critério: "  [plain]value+br
#   [bold]"Norma de recuperação: "  [plain]value
# Work from the last run to the first so earlier anchors stay valid.
$pAval = $d.Paragraphs.Item(14).Range

# "Norma de recuperação:" value -> last run, no trailing break
$avalText = $pAval.Text
$normaLabel = "Norma de recupera" + [char]0xE7 + [char]0xE3 + "o: "
$normaIdx = $avalText.IndexOf($normaLabel)
if ($normaIdx -lt 0) {
    Write-Output "WARNING: 'Norma de recuperacao:' label not found"
}
$s = $pAval.Start + $normaIdx + $normaLabel.Length
$e = $pAval.End - 1
$r = $d.Range($s, $e)
$r.Text = "471420 - Carlos Antonio Reis Pereira Baptista"

# "Critério:" value -> now split into two text segments joined by a
# <w:br/>, still followed by the break that was already there.
$avalText = $pAval.Text
$criterioLabel = "Crit" + [char]0xE9 + "rio: "
$criterioIdx = $avalText.IndexOf($criterioLabel)
if ($criterioIdx -lt 0) {
    Write-Output "WARNING: 'Criterio:' label not found"
}
$criterioValStart = $pAval.Start + $criterioIdx + $criterioLabel.Length
$criterioBrIdx = $avalText.IndexOf($vtab, $criterioIdx)
$s = $criterioValStart
$e = $pAval.Start + $criterioBrIdx + 1
$r = $d.Range($s, $e)
$newCriterio = "Apostila ou texto fornecido pelo(s) docente(s) responsáveis." + $vtab + "Artigos extraídos de revistas especializadas na área de Ciência e Engenharia de Materiais." + $vtab
$r.Text = $newCriterio

# "Método:" value -> first run, keeps its trailing break
$avalText = $pAval.Text
$metodoLabel = "M" + [char]0xE9 + "todo: "
$metodoIdx = $avalText.IndexOf($metodoLabel)
if ($metodoIdx -lt 0) {
    Write-Output "WARNING: 'Metodo:' label not found"
}
$metodoValStart = $pAval.Start + $metodoIdx + $metodoLabel.Length
$metodoBrIdx = $avalText.IndexOf($vtab)
$s = $metodoValStart
$e = $pAval.Start + $metodoBrIdx
$r = $d.Range($s, $e)
$r.Text = "Não cabe recuperação."

# --- "Bibliografia" body (paragraph 16): collapses two runs + <w:br/>
#     into a single run ---
$pBiblio = $d.Paragraphs.Item(16).Range
$s = $pBiblio.Start
$e = $pBiblio.End - 1
$r = $d.Range($s, $e)
$r.Text = "519033 - Carlos Yujiro Shigue"

Write-Output "done"
